$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$lo = $ws.ListObjects.Item(1)
$lo.ListRows.Add() | Out-Null
$lo.ListRows.Add() | Out-Null

$ws.Range("A2").Value = "coin01"
$ws.Range("B2").Value = "Sprites/GoldenCoin"
$ws.Range("C2").Value = "Coin"
$ws.Range("D2").Value = "A Golden Coin"

$ws.Range("A3").Value = "star01"
$ws.Range("B3").Value = "Sprites/star"
$ws.Range("C3").Value = "Star"
$ws.Range("D3").Value = "A Golden Star"
$ws.Range("A3:D3").Style = "Normal"

$ws.Range("A4").Value = "shuffle01"
$ws.Range("B4").Value = "Sprites/shuffle"
$ws.Range("C4").Value = "Shuffle"
$ws.Range("D4").Value = "Shuffle Icon"

$ws.Range("A5").Value = "circle"
$ws.Range("B5").Value = "Sprites/Circle"
$ws.Range("C5").Value = "Circle"
$ws.Range("D5").Value = "A Circle"

$ws.Range("C3").Select() | Out-Null
